$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Daily price-data refresh: insert a new top row for 2025-12-17, pushing all
# existing date rows down by one (row 2 -> 3, ... , row 27 -> 28).
$ws.Rows.Item(2).Insert()

# Format the date cell as text first so Excel stores the literal
# "2025-12-17" string (matching the other date cells) instead of
# auto-converting it to a date serial value.
$ws.Cells.Item(2, 1).NumberFormat = "@"
$ws.Cells.Item(2, 1).Value = "2025-12-17"

# New day's values match the prior latest day's readings.
$ws.Cells.Item(2, 2).Value = 783.5
$ws.Cells.Item(2, 3).Value = 1112
$ws.Cells.Item(2, 4).Value = 3610

# Drop the formatting picked up from the inserted row (and the temporary
# text number format) so the new row matches the plain, unstyled data rows.
$ws.Rows.Item(2).ClearFormats()
